$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 20 with forecast data, continuing the existing series pattern.
# Copy the formatting of the A column date cell above (row 19) down to row 20
# so the new date cell keeps the same style (bold, border, centered, date format).
$ws.Range("A19").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = -0.08656168856399082
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = -0.1516437243033186
